# Weekly update: insert a new "Arveja Verde" record for Vega Modelo de Temuco
# at row 55, shifting the existing rows 55-88 down to 56-89.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 55; everything currently at 55..88 moves to 56..89.
$ws.Rows.Item(55).Insert()

# Populate the newly inserted row 55 with the new weekly record.
$ws.Cells.Item(55, 1).Value = 10
$ws.Cells.Item(55, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(55, 3).Value = "La Araucanía"
$ws.Cells.Item(55, 4).Value = 44574
$ws.Cells.Item(55, 5).Value = 9
$ws.Cells.Item(55, 6).Value = 100112022
$ws.Cells.Item(55, 7).Value = "Arveja Verde"
$ws.Cells.Item(55, 8).Value = "Sin especificar"
$ws.Cells.Item(55, 9).Value = "Primera"
$ws.Cells.Item(55, 10).Value = 235
$ws.Cells.Item(55, 11).Value = 25000
$ws.Cells.Item(55, 12).Value = 27000
$ws.Cells.Item(55, 13).Value = 26064
$ws.Cells.Item(55, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(55, 15).Value = "Región de La Araucanía"
$ws.Cells.Item(55, 16).Value = 1043
$ws.Cells.Item(55, 17).Value = 25
$ws.Cells.Item(55, 18).Value = "Hortaliza"

# Keep the date cell formatted the same way as the other date cells in column D.
$ws.Cells.Item(55, 4).NumberFormat = $ws.Cells.Item(56, 4).NumberFormat
